$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.479.25'
$ws.Range('E2').Value = '  -2.63%  '
$ws.Range('D3').Value = '1.776.77'
$ws.Range('E3').Value = '  -2.92%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''229.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '''0.5871'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('D8').Value = '''0.2744'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').Value = '''23.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '''0.06691'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('D11').Value = '''0.07527'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').Value = '1.792.64'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').Value = '''4.754'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '''0.6084'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.24%  '
$ws.Range('D15').Value = '2.018.08'
$ws.Range('E15').Value = '  -3.00%  '
$ws.Range('D16').Value = '''74.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.67%  '
$ws.Range('D17').Value = '''0.000008647'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -11.04%  '
$ws.Range('D18').Value = '28.439.86'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('E19').Value = '  -5.54%  '
$ws.Range('D21').Value = '''207.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.32%  '
$ws.Range('D22').Value = '''11.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').Value = '''6.739'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '''151.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.06%  '
$ws.Range('D26').Value = '''8.109'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').Value = '''0.1247'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('D28').Value = '''16.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').Value = '''1.410'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('D30').Value = '''0.06156'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.11%  '
$ws.Range('D31').Value = '''1.414'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('D32').Value = '''3.762'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').Value = '''3.751'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').Value = '''1.674'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('D35').Value = '''1.044'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.49%  '
$ws.Range('D36').Value = '''0.6361'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('D38').Value = '''2.681'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').Value = '1.144.09'
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = '''0.01674'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('D41').Value = '''6.277'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.79%  '
$ws.Range('D42').Value = '''0.8737'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('D43').Value = '''1.006'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '''99.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = '1.929.19'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('D46').Value = '''59.63'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').Value = '''0.00000000109'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.84%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.576'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''8.389'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('D50').Value = '''0.05414'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('D51').Value = '''0.4463'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.97%  '
